$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1061.8
$ws.Range("I19").Value = 1205.8334
$ws.Range("J19").Value = 965.7778
$ws.Range("K19").Value = 1205.8334
$ws.Range("L19").Value = 965.7778
$ws.Range("M19").Value = -1030.8334
$ws.Range("N19").Value = -1315.7778
$ws.Range("H132").Value = 12236.946
$ws.Range("I132").Value = 2802.375
$ws.Range("K132").Value = 8407.125
$ws.Range("M132").Value = -5877.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5757.485
$ws.Range("I32").Value = 5894.6772
$ws.Range("K32").Value = 5894.6772
$ws.Range("M32").Value = -5607.6772
$ws.Range("H61").Value = 2844.4375
$ws.Range("I61").Value = 1069.1111
$ws.Range("K61").Value = 1069.1111
$ws.Range("M61").Value = -857.1111000000001
$ws.Range("H74").Value = 2233.6667
$ws.Range("I74").Value = 1380.625
$ws.Range("J74").Value = 3939.75
$ws.Range("K74").Value = 1380.625
$ws.Range("L74").Value = 3939.75
$ws.Range("M74").Value = -506.625
$ws.Range("N74").Value = -5687.75
$ws.Range("H77").Value = 2233.6667
$ws.Range("I77").Value = 1380.625
$ws.Range("J77").Value = 3939.75
$ws.Range("K77").Value = 6903.125
$ws.Range("L77").Value = 19698.75
$ws.Range("M77").Value = -2535.125
$ws.Range("N77").Value = -28434.75
$ws.Range("H97").Value = 5558405.5
$ws.Range("I97").Value = 2118.1177
$ws.Range("J97").Value = 37044030
$ws.Range("K97").Value = 2118.1177
$ws.Range("L97").Value = 37044030
$ws.Range("M97").Value = -1622.1177
$ws.Range("N97").Value = -37045022
$ws.Range("H122").Value = 3319.15
$ws.Range("I122").Value = 2724
$ws.Range("J122").Value = 4211.875
$ws.Range("K122").Value = 8172
$ws.Range("L122").Value = 12635.625
$ws.Range("M122").Value = -5722
$ws.Range("N122").Value = -17535.625
$ws.Range("H132").Value = 1466.1333
$ws.Range("I132").Value = 1071.1428
$ws.Range("J132").Value = 1811.75
$ws.Range("K132").Value = 3213.4284
$ws.Range("L132").Value = 5435.25
$ws.Range("M132").Value = -683.4284000000002
$ws.Range("N132").Value = -10495.25
$ws.Range("H136").Value = 2844.4375
$ws.Range("I136").Value = 1069.1111
$ws.Range("K136").Value = 3207.3333
$ws.Range("M136").Value = -657.3333000000002

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 16670826
$ws.Range("I94").Value = 4620.3335
$ws.Range("J94").Value = 166666670
$ws.Range("K94").Value = 4620.3335
$ws.Range("L94").Value = 166666670
$ws.Range("M94").Value = -4169.3335
$ws.Range("N94").Value = -166667572
$ws.Range("H95").Value = 4500
$ws.Range("J95").Value = 4500
$ws.Range("L95").Value = 4500
$ws.Range("N95").Value = -9992
$ws.Range("H97").Value = 10183.308
$ws.Range("I97").Value = 8264.727999999999
$ws.Range("J97").Value = 20735.5
$ws.Range("K97").Value = 8264.727999999999
$ws.Range("L97").Value = 20735.5
$ws.Range("M97").Value = -7273.727999999999
$ws.Range("N97").Value = -22717.5
$ws.Range("H99").Value = 19723.61
$ws.Range("I99").Value = 25575.117
$ws.Range("J99").Value = 3144.3333
$ws.Range("K99").Value = 25575.117
$ws.Range("L99").Value = 3144.3333
$ws.Range("M99").Value = -24077.117
$ws.Range("N99").Value = -6140.3333
$ws.Range("H107").Value = 4893.2666
$ws.Range("I107").Value = 3814.2144
$ws.Range("K107").Value = 3814.2144
$ws.Range("M107").Value = -1894.2144

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 505.2
$ws.Range("I22").Value = 371.875
$ws.Range("J22").Value = 1038.5
$ws.Range("K22").Value = 371.875
$ws.Range("L22").Value = 1038.5
$ws.Range("M22").Value = -21.875
$ws.Range("N22").Value = -1738.5
$ws.Range("H31").Value = 2447.7
$ws.Range("I31").Value = 2347.3076
$ws.Range("J31").Value = 2634.1428
$ws.Range("K31").Value = 2347.3076
$ws.Range("L31").Value = 2634.1428
$ws.Range("M31").Value = -2052.3076
$ws.Range("N31").Value = -3224.1428
$ws.Range("H34").Value = 2447.7
$ws.Range("I34").Value = 2347.3076
$ws.Range("J34").Value = 2634.1428
$ws.Range("K34").Value = 2347.3076
$ws.Range("L34").Value = 2634.1428
$ws.Range("M34").Value = -2145.3076
$ws.Range("N34").Value = -3038.1428
$ws.Range("H53").Value = 95166
$ws.Range("J53").Value = 95166
$ws.Range("L53").Value = 95166
$ws.Range("N53").Value = -96380
$ws.Range("H86").Value = 40548.332
$ws.Range("I86").Value = 45251.375
$ws.Range("J86").Value = 31142.25
$ws.Range("K86").Value = 45251.375
$ws.Range("L86").Value = 31142.25
$ws.Range("M86").Value = -44128.375
$ws.Range("N86").Value = -33388.25
$ws.Range("H89").Value = 40548.332
$ws.Range("I89").Value = 45251.375
$ws.Range("J89").Value = 31142.25
$ws.Range("K89").Value = 226256.875
$ws.Range("L89").Value = 155711.25
$ws.Range("M89").Value = -220640.875
$ws.Range("N89").Value = -166943.25
$ws.Range("H105").Value = 3887.6667
$ws.Range("I105").Value = 772
$ws.Range("K105").Value = 772
$ws.Range("M105").Value = 975
$ws.Range("H107").Value = 11799.368
$ws.Range("I107").Value = 1316.3334
$ws.Range("K107").Value = 1316.3334
$ws.Range("M107").Value = 603.6666
$ws.Range("H122").Value = 412342.2
$ws.Range("I122").Value = 681637.25
$ws.Range("K122").Value = 2044911.75
$ws.Range("M122").Value = -2042461.75
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1879.88
$ws.Range("I132").Value = 1933.5814
$ws.Range("K132").Value = 5800.7442
$ws.Range("M132").Value = -3270.7442
$ws.Range("H141").Value = 46666.25
$ws.Range("J141").Value = 46666.25
$ws.Range("L141").Value = 46666.25
$ws.Range("N141").Value = -57026.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 62503350
$ws.Range("I139").Value = 90910060
$ws.Range("K139").Value = 272730180
$ws.Range("M139").Value = -272725040

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4594.184
$ws.Range("I122").Value = 3799.3157
$ws.Range("J122").Value = 5389.0527
$ws.Range("K122").Value = 11397.9471
$ws.Range("L122").Value = 16167.1581
$ws.Range("M122").Value = -8947.947100000001
$ws.Range("N122").Value = -21067.1581
$ws.Range("H132").Value = 2944.75
$ws.Range("I132").Value = 3174.7778
$ws.Range("J132").Value = 2254.6667
$ws.Range("K132").Value = 9524.3334
$ws.Range("L132").Value = 6764.000100000001
$ws.Range("M132").Value = -6994.3334
$ws.Range("N132").Value = -11824.0001
$ws.Range("H137").Value = 70780
$ws.Range("J137").Value = 70780
$ws.Range("L137").Value = 70780
$ws.Range("N137").Value = -80980

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 761.9167
$ws.Range("I16").Value = 735.3
$ws.Range("K16").Value = 735.3
$ws.Range("M16").Value = -565.3
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H46").Value = 9354.454
$ws.Range("I46").Value = 3749.75
$ws.Range("K46").Value = 3749.75
$ws.Range("M46").Value = -3561.75
$ws.Range("H55").Value = 2754.7222
$ws.Range("J55").Value = 3150
$ws.Range("L55").Value = 3150
$ws.Range("N55").Value = -3496
$ws.Range("H122").Value = 7032
$ws.Range("I122").Value = 2899.5
$ws.Range("K122").Value = 8698.5
$ws.Range("M122").Value = -6248.5
$ws.Range("H132").Value = 3529.5425
$ws.Range("I132").Value = 2299.7551
$ws.Range("K132").Value = 6899.265299999999
$ws.Range("M132").Value = -4369.265299999999
$ws.Range("H141").Value = 77000
$ws.Range("J141").Value = 74833.336
$ws.Range("L141").Value = 74833.336
$ws.Range("N141").Value = -85193.336

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 283451.4
$ws.Range("I122").Value = 2451.2334
$ws.Range("J122").Value = 1126451.9
$ws.Range("K122").Value = 7353.7002
$ws.Range("L122").Value = 3379355.7
$ws.Range("M122").Value = -4903.7002
$ws.Range("N122").Value = -3384255.7
$ws.Range("H132").Value = 1653.9395
$ws.Range("I132").Value = 1316.5862
$ws.Range("K132").Value = 3949.7586
$ws.Range("M132").Value = -1419.7586
$ws.Range("H136").Value = 1773.9445
$ws.Range("I136").Value = 1095.0714
$ws.Range("K136").Value = 3285.2142
$ws.Range("M136").Value = -735.2142000000003
